$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.028408779208068
$ws.Range("D2").Value = 1.030595463669476
$ws.Range("E2").Value = 1.037033332487
$ws.Range("F2").Value = 1.044843750138388
$ws.Range("I2").Value = 1.029700318295664
$ws.Range("J2").Value = 1.033561260467731
$ws.Range("K2").Value = 1.033406019534925
$ws.Range("L2").Value = 1.039825372727055
$ws.Range("M2").Value = 1.047613659673579

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.029391606141853
$ws.Range("D3").Value = 1.031280565390298
$ws.Range("E3").Value = 1.037942169326056
$ws.Range("F3").Value = 1.045919146538698
$ws.Range("I3").Value = 1.029818778292732
$ws.Range("J3").Value = 1.034184513075173
$ws.Range("K3").Value = 1.03389984014603
$ws.Range("L3").Value = 1.040543671192461
$ws.Range("M3").Value = 1.04849968785429

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.03002794556716
$ws.Range("D4").Value = 1.031723906649634
$ws.Range("E4").Value = 1.038531017501686
$ws.Range("F4").Value = 1.046616102736138
$ws.Range("I4").Value = 1.029894055214282
$ws.Range("J4").Value = 1.034587598436207
$ws.Range("K4").Value = 1.034218733465333
$ws.Range("L4").Value = 1.041008598321494
$ws.Range("M4").Value = 1.049073513562051

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.030295554107258
$ws.Range("D5").Value = 1.031910294321527
$ws.Range("E5").Value = 1.038778752586523
$ws.Range("F5").Value = 1.046909366251542
$ws.Range("I5").Value = 1.029925372055578
$ws.Range("J5").Value = 1.034757006672969
$ws.Range("K5").Value = 1.03435264170406
$ws.Range("L5").Value = 1.041204086374852
$ws.Range("M5").Value = 1.049314870233168

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.030340492102388
$ws.Range("D6").Value = 1.031941590005252
$ws.Range("E6").Value = 1.038820359117235
$ws.Range("F6").Value = 1.046958621890598
$ws.Range("I6").Value = 1.029930610950264
$ws.Range("J6").Value = 1.03478544817393
$ws.Range("K6").Value = 1.034375116399599
$ws.Range("L6").Value = 1.041236911572099
$ws.Range("M6").Value = 1.049355402112687

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.030031521007052
$ws.Range("D7").Value = 1.031726397144918
$ws.Range("E7").Value = 1.038534327031002
$ws.Range("F7").Value = 1.046620020304768
$ws.Range("I7").Value = 1.029894474967072
$ws.Range("J7").Value = 1.03458986227029
$ws.Range("K7").Value = 1.03422052336204
$ws.Range("L7").Value = 1.041011210315022
$ws.Range("M7").Value = 1.049076738109896

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.028740850130886
$ws.Range("D8").Value = 1.030826989060154
$ws.Range("E8").Value = 1.037340318110652
$ws.Range("F8").Value = 1.045206956506708
$ws.Range("I8").Value = 1.029740636735092
$ws.Range("J8").Value = 1.033771932696776
$ws.Range("K8").Value = 1.033573040834418
$ws.Range("L8").Value = 1.040068095573263
$ws.Range("M8").Value = 1.047912991716348

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.02646949779077
$ws.Range("D9").Value = 1.02924243896241
$ws.Range("E9").Value = 1.035242263834506
$ws.Range("F9").Value = 1.042725442874358
$ws.Range("I9").Value = 1.02945904779942
$ws.Range("J9").Value = 1.032329133485793
$ws.Range("K9").Value = 1.032427219506918
$ws.Range("L9").Value = 1.038407315427914
$ws.Range("M9").Value = 1.045866236354213

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.024957299426451
$ws.Range("D10").Value = 1.0281863641626
$ws.Range("E10").Value = 1.033847615444408
$ws.Range("F10").Value = 1.041076856586623
$ws.Range("I10").Value = 1.029264286794394
$ws.Range("J10").Value = 1.031366302726626
$ws.Range("K10").Value = 1.031660114215001
$ws.Range("L10").Value = 1.0373009219029
$ws.Range("M10").Value = 1.044504413679536

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.024302990978667
$ws.Range("D11").Value = 1.027729157034499
$ws.Range("E11").Value = 1.033244691091344
$ws.Range("F11").Value = 1.04036437605686
$ws.Range("I11").Value = 1.029178290457414
$ws.Range("J11").Value = 1.030949166894842
$ws.Range("K11").Value = 1.031327195251961
$ws.Range("L11").Value = 1.036822039180131
$ws.Range("M11").Value = 1.043915375194165

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.0240600248135
$ws.Range("D12").Value = 1.027559343236209
$ws.Range("E12").Value = 1.033020884446943
$ws.Range("F12").Value = 1.040099935392317
$ws.Range("I12").Value = 1.029146098109302
$ws.Range("J12").Value = 1.030794191160425
$ws.Range("K12").Value = 1.031203421385566
$ws.Range("L12").Value = 1.036644190545716
$ws.Range("M12").Value = 1.043696676927822

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.024112138563593
$ws.Range("D13").Value = 1.02759576825139
$ws.Range("E13").Value = 1.033068885096248
$ws.Range("F13").Value = 1.040156649463388
$ws.Range("I13").Value = 1.02915301475953
$ws.Range("J13").Value = 1.030827435459628
$ws.Range("K13").Value = 1.03122997639742
$ws.Range("L13").Value = 1.036682338324245
$ws.Range("M13").Value = 1.043743584053522

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.024282905825033
$ws.Range("D14").Value = 1.027715119894802
$ws.Range("E14").Value = 1.033226188174304
$ws.Range("F14").Value = 1.040342513095556
$ws.Range("I14").Value = 1.029175634514591
$ws.Range("J14").Value = 1.03093635721964
$ws.Range("K14").Value = 1.031316966364662
$ws.Range("L14").Value = 1.036807337545746
$ws.Range("M14").Value = 1.043897295547141

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.024388130875612
$ws.Range("D15").Value = 1.0277886581788
$ws.Range("E15").Value = 1.033323127207227
$ws.Range("F15").Value = 1.040457057177072
$ws.Range("I15").Value = 1.029189538247178
$ws.Range("J15").Value = 1.031003463170573
$ws.Range("K15").Value = 1.031370548813855
$ws.Range("L15").Value = 1.036884357653397
$ws.Range("M15").Value = 1.043992015127879

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.025000733950656
$ws.Range("D16").Value = 1.028216709289526
$ws.Range("E16").Value = 1.033887650014017
$ws.Range("F16").Value = 1.041124170465087
$ws.Range("I16").Value = 1.029269959077481
$ws.Range("J16").Value = 1.031393982003669
$ws.Range("K16").Value = 1.031682193065935
$ws.Range("L16").Value = 1.037332707900824
$ws.Range("M16").Value = 1.044543519765398

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.025385133450057
$ws.Range("D17").Value = 1.028485236792924
$ws.Range("E17").Value = 1.03424201997003
$ws.Range("F17").Value = 1.041542999856663
$ws.Range("I17").Value = 1.029319959879769
$ws.Range("J17").Value = 1.031638884740069
$ws.Range("K17").Value = 1.031877477109406
$ws.Range("L17").Value = 1.037613998351641
$ws.Range("M17").Value = 1.044889635848961

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.025609393784273
$ws.Range("D18").Value = 1.028641872000938
$ws.Range("E18").Value = 1.034448811266848
$ws.Range("F18").Value = 1.041787427915551
$ws.Range("I18").Value = 1.029348963964565
$ws.Range("J18").Value = 1.031781710725981
$ws.Range("K18").Value = 1.031991309862515
$ws.Range("L18").Value = 1.037778088882155
$ws.Range("M18").Value = 1.045091581178302

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.025685868653618
$ws.Range("D19").Value = 1.02869528182052
$ws.Range("E19").Value = 1.034519337551482
$ws.Range("F19").Value = 1.041870793926549
$ws.Range("I19").Value = 1.029358826364427
$ws.Range("J19").Value = 1.031830407004855
$ws.Range("K19").Value = 1.032030111460807
$ws.Range("L19").Value = 1.037834042639858
$ws.Range("M19").Value = 1.045160449755217

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.025343886196973
$ws.Range("D20").Value = 1.028456425532049
$ws.Range("E20").Value = 1.034203989770542
$ws.Range("F20").Value = 1.04149804977177
$ws.Range("I20").Value = 1.029314611868759
$ws.Range("J20").Value = 1.031612611227063
$ws.Range("K20").Value = 1.031856532537329
$ws.Range("L20").Value = 1.037583816622255
$ws.Range("M20").Value = 1.044852494484191

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.024232617100169
$ws.Range("D21").Value = 1.027679973433924
$ws.Range("E21").Value = 1.033179862292234
$ws.Range("F21").Value = 1.04028777518657
$ws.Range("I21").Value = 1.029168980443765
$ws.Range("J21").Value = 1.030904283379036
$ws.Range("K21").Value = 1.031291353098463
$ws.Range("L21").Value = 1.036770527563872
$ws.Range("M21").Value = 1.043852028663797

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.02353434144357
$ws.Range("D22").Value = 1.027191864661261
$ws.Range("E22").Value = 1.032536799641506
$ws.Range("F22").Value = 1.039528021649661
$ws.Range("I22").Value = 1.0290759729971
$ws.Range("J22").Value = 1.03045873946947
$ws.Range("K22").Value = 1.030935349329697
$ws.Range("L22").Value = 1.03625935283546
$ws.Range("M22").Value = 1.043223556832239

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.023904470286971
$ws.Range("D23").Value = 1.027450612601855
$ws.Range("E23").Value = 1.032877618727109
$ws.Range("F23").Value = 1.039930667926974
$ws.Range("I23").Value = 1.029125414652312
$ws.Range("J23").Value = 1.030694948507764
$ws.Range("K23").Value = 1.031124135283062
$ws.Range("L23").Value = 1.036530319643754
$ws.Range("M23").Value = 1.04355666824068

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.025362523911122
$ws.Range("D24").Value = 1.028469444077118
$ws.Range("E24").Value = 1.03422117369248
$ws.Range("F24").Value = 1.041518360372554
$ws.Range("I24").Value = 1.029317028900798
$ws.Range("J24").Value = 1.031624483163542
$ws.Range("K24").Value = 1.031865996713903
$ws.Range("L24").Value = 1.03759745438926
$ws.Range("M24").Value = 1.044869276877655

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.027056340554026
$ws.Range("D25").Value = 1.029652037210649
$ws.Range("E25").Value = 1.035783950882585
$ws.Range("F25").Value = 1.043365962825395
$ws.Range("I25").Value = 1.029533087374768
$ws.Range("J25").Value = 1.032702305041386
$ws.Range("K25").Value = 1.032724013977079
$ws.Range("L25").Value = 1.038836530604861
$ws.Range("M25").Value = 1.046394902748659
